# Applies the updated currentAveragePrice / LevePrice / LeveProfit
# figures (columns H-N) captured by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 12665.5
$ws.Range("I28").Value = 12198.9
$ws.Range("K28").Value = 12198.9
$ws.Range("M28").Value = -11713.9
$ws.Range("H29").Value = 3323.75
$ws.Range("I29").Value = 648.5
$ws.Range("K29").Value = 1945.5
$ws.Range("M29").Value = -1664.5
$ws.Range("H38").Value = 825.8333
$ws.Range("J38").Value = 4000
$ws.Range("L38").Value = 12000
$ws.Range("N38").Value = -12744
$ws.Range("H58").Value = 646.5
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("H64").Value = 6978.6
$ws.Range("J64").Value = 7623.25
$ws.Range("L64").Value = 7623.25
$ws.Range("N64").Value = -8119.25
$ws.Range("H67").Value = 6978.6
$ws.Range("J67").Value = 7623.25
$ws.Range("L67").Value = 7623.25
$ws.Range("N67").Value = -9339.25
$ws.Range("H74").Value = 6222.853
$ws.Range("I74").Value = 6282.6875
$ws.Range("K74").Value = 6282.6875
$ws.Range("M74").Value = -5346.6875
$ws.Range("H77").Value = 6222.853
$ws.Range("I77").Value = 6282.6875
$ws.Range("K77").Value = 31413.4375
$ws.Range("M77").Value = -26733.4375
$ws.Range("H98").Value = 1801.7097
$ws.Range("I98").Value = 1830.5333
$ws.Range("K98").Value = 1830.5333
$ws.Range("M98").Value = -332.5333000000001
$ws.Range("H122").Value = 1801.7097
$ws.Range("I122").Value = 1830.5333
$ws.Range("K122").Value = 5491.5999
$ws.Range("M122").Value = -3041.5999
$ws.Range("H134").Value = 40530.266
$ws.Range("J134").Value = 38068.145
$ws.Range("L134").Value = 38068.145
$ws.Range("N134").Value = -48208.145
$ws.Range("H137").Value = 2464.1667
$ws.Range("I137").Value = 2457
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 7371
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -4821
$ws.Range("N137").Value = -12600
$ws.Range("H138").Value = 8045.9736
$ws.Range("I138").Value = 14384.714
$ws.Range("J138").Value = 6614.645
$ws.Range("K138").Value = 43154.142
$ws.Range("L138").Value = 19843.935
$ws.Range("M138").Value = -38014.142
$ws.Range("N138").Value = -30123.935
$ws.Range("N58").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26591.666
$ws.Range("I32").Value = 19217.857
$ws.Range("K32").Value = 19217.857
$ws.Range("M32").Value = -18930.857
$ws.Range("H61").Value = 4171.4688
$ws.Range("I61").Value = 2782.7778
$ws.Range("J61").Value = 11670.4
$ws.Range("K61").Value = 2782.7778
$ws.Range("L61").Value = 11670.4
$ws.Range("M61").Value = -2570.7778
$ws.Range("N61").Value = -12094.4
$ws.Range("H122").Value = 9738.799999999999
$ws.Range("I122").Value = 9738.799999999999
$ws.Range("K122").Value = 29216.4
$ws.Range("M122").Value = -26766.4
$ws.Range("H132").Value = 18185504
$ws.Range("I132").Value = 25002930
$ws.Range("J132").Value = 5700.533
$ws.Range("K132").Value = 75008790
$ws.Range("L132").Value = 17101.599
$ws.Range("M132").Value = -75006260
$ws.Range("N132").Value = -22161.599
$ws.Range("H136").Value = 4171.4688
$ws.Range("I136").Value = 2782.7778
$ws.Range("J136").Value = 11670.4
$ws.Range("K136").Value = 8348.3334
$ws.Range("L136").Value = 35011.2
$ws.Range("M136").Value = -5798.3334
$ws.Range("N136").Value = -40111.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 44399
$ws.Range("J60").Value = 44399
$ws.Range("L60").Value = 44399
$ws.Range("N60").Value = -45597
$ws.Range("H64").Value = 1727
$ws.Range("I64").Value = 1283.8
$ws.Range("J64").Value = 2465.6667
$ws.Range("K64").Value = 1283.8
$ws.Range("L64").Value = 2465.6667
$ws.Range("M64").Value = -1058.8
$ws.Range("N64").Value = -2915.6667
$ws.Range("H67").Value = 1727
$ws.Range("I67").Value = 1283.8
$ws.Range("J67").Value = 2465.6667
$ws.Range("K67").Value = 1283.8
$ws.Range("L67").Value = 2465.6667
$ws.Range("M67").Value = -503.8
$ws.Range("N67").Value = -4025.6667
$ws.Range("H94").Value = 26418.25
$ws.Range("I94").Value = 1891
$ws.Range("J94").Value = 100000
$ws.Range("K94").Value = 1891
$ws.Range("L94").Value = 100000
$ws.Range("M94").Value = -1440
$ws.Range("N94").Value = -100902
$ws.Range("H107").Value = 980.5
$ws.Range("I107").Value = 979.6667
$ws.Range("K107").Value = 979.6667
$ws.Range("M107").Value = 940.3333
$ws.Range("H128").Value = 1000
$ws.Range("I128").Value = 1000
$ws.Range("K128").Value = 3000
$ws.Range("M128").Value = -510
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 979
$ws.Range("I16").Value = 979
$ws.Range("K16").Value = 979
$ws.Range("M16").Value = -692
$ws.Range("H31").Value = 146555.72
$ws.Range("I31").Value = 4231.6665
$ws.Range("K31").Value = 4231.6665
$ws.Range("M31").Value = -3936.6665
$ws.Range("H34").Value = 146555.72
$ws.Range("I34").Value = 4231.6665
$ws.Range("K34").Value = 4231.6665
$ws.Range("M34").Value = -4029.6665
$ws.Range("H107").Value = 368.85715
$ws.Range("I107").Value = 326.93332
$ws.Range("K107").Value = 326.93332
$ws.Range("M107").Value = 1593.06668
$ws.Range("H113").Value = 979
$ws.Range("I113").Value = 979
$ws.Range("K113").Value = 979
$ws.Range("M113").Value = 1191
$ws.Range("H132").Value = 100005704
$ws.Range("I132").Value = 200004800
$ws.Range("J132").Value = 6600.6
$ws.Range("K132").Value = 600014400
$ws.Range("L132").Value = 19801.8
$ws.Range("M132").Value = -600011870
$ws.Range("N132").Value = -24861.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7064.08
$ws.Range("I56").Value = 7064.08
$ws.Range("K56").Value = 7064.08
$ws.Range("M56").Value = -6534.08
$ws.Range("H103").Value = 1166
$ws.Range("I103").Value = 1166
$ws.Range("K103").Value = 3498
$ws.Range("M103").Value = -2619
$ws.Range("H114").Value = 2722.8
$ws.Range("I114").Value = 566.6667
$ws.Range("J114").Value = 3646.8572
$ws.Range("K114").Value = 1700.0001
$ws.Range("L114").Value = 10940.5716
$ws.Range("M114").Value = 1553.9999
$ws.Range("N114").Value = -17448.5716
$ws.Range("H140").Value = 6429674.5
$ws.Range("I140").Value = 16669155
$ws.Range("J140").Value = 29999.166
$ws.Range("K140").Value = 50007465
$ws.Range("L140").Value = 89997.49800000001
$ws.Range("M140").Value = -50002285
$ws.Range("N140").Value = -100357.498

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7481.48
$ws.Range("I70").Value = 7545.1177
$ws.Range("J70").Value = 7346.25
$ws.Range("K70").Value = 7545.1177
$ws.Range("L70").Value = 7346.25
$ws.Range("M70").Value = -7275.1177
$ws.Range("N70").Value = -7886.25
$ws.Range("H73").Value = 7481.48
$ws.Range("I73").Value = 7545.1177
$ws.Range("J73").Value = 7346.25
$ws.Range("K73").Value = 7545.1177
$ws.Range("L73").Value = 7346.25
$ws.Range("M73").Value = -6609.1177
$ws.Range("N73").Value = -9218.25
$ws.Range("H99").Value = 31799.375
$ws.Range("I99").Value = 20939.6
$ws.Range("J99").Value = 49899
$ws.Range("K99").Value = 20939.6
$ws.Range("L99").Value = 49899
$ws.Range("M99").Value = -18693.6
$ws.Range("N99").Value = -54391
$ws.Range("H113").Value = 3580.5417
$ws.Range("I113").Value = 3338
$ws.Range("J113").Value = 3867.182
$ws.Range("K113").Value = 3338
$ws.Range("L113").Value = 3867.182
$ws.Range("M113").Value = -1168
$ws.Range("N113").Value = -8207.182000000001
$ws.Range("H132").Value = 20018410
$ws.Range("I132").Value = 50006704
$ws.Range("K132").Value = 150020112
$ws.Range("M132").Value = -150017582

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1081.6818
$ws.Range("I22").Value = 889.6667
$ws.Range("K22").Value = 889.6667
$ws.Range("M22").Value = -594.6667
$ws.Range("H27").Value = 1081.6818
$ws.Range("I27").Value = 889.6667
$ws.Range("K27").Value = 889.6667
$ws.Range("M27").Value = -782.6667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10544.5
$ws.Range("J81").Value = 4600
$ws.Range("L81").Value = 9200
$ws.Range("N81").Value = -11322
$ws.Range("H84").Value = 10544.5
$ws.Range("J84").Value = 4600
$ws.Range("L84").Value = 46000
$ws.Range("N84").Value = -56608
$ws.Range("H126").Value = 7146037
$ws.Range("I126").Value = 8623893
$ws.Range("J126").Value = 3066.1667
$ws.Range("K126").Value = 25871679
$ws.Range("L126").Value = 9198.500100000001
$ws.Range("M126").Value = -25869209
$ws.Range("N126").Value = -14138.5001
